# Refresh the scraped crypto-price table (cryptos.xlsx): on each GitHub Actions
# run the scraper re-pulls Price (column D) and Volume(1h) (column E) for every
# coin row and overwrites the corresponding cells.
#
# The existing cells hold plain text (no number format applied), e.g. D2 is the
# literal string "261.59" and E2 is the literal string "0.89%" - not real numbers
# or percentages. Writing straight to `.Value`/`.Value2`/`.Formula` with a
# numeric- or percent-looking string triggers Excel's automatic type detection,
# which would silently turn the text into a real number/percentage and assign a
# new number format to the cell - changing more than just the displayed text.
#
# To keep the edit faithful to a plain-text overwrite, each cell is first given a
# formula that evaluates to the desired text (quoted, so the formula result is a
# string), then immediately copy/paste-special'd (values only) onto itself. That
# commits the literal text as the stored value while leaving the cell's original
# (unformatted) style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value="261.73"},
    @{Cell="E2"; Value="0.94%"},
    @{Cell="D3"; Value="27.12"},
    @{Cell="E3"; Value="0.34%"},
    @{Cell="D4"; Value="4.710"},
    @{Cell="E4"; Value="0.09%"},
    @{Cell="D5"; Value="0.06206"},
    @{Cell="E5"; Value="2.93%"},
    @{Cell="E6"; Value="0.80%"},
    @{Cell="D7"; Value="0.8501"},
    @{Cell="E7"; Value="-1.16%"},
    @{Cell="D8"; Value="0.9112"},
    @{Cell="E8"; Value="-1.42%"},
    @{Cell="D9"; Value="0.1408"},
    @{Cell="E9"; Value="0.88%"},
    @{Cell="D10"; Value="0.04694"},
    @{Cell="E10"; Value="-7.87%"},
    @{Cell="D11"; Value="0.07087"},
    @{Cell="E11"; Value="0.10%"},
    @{Cell="D12"; Value="0.03158"},
    @{Cell="E12"; Value="2.38%"},
    @{Cell="D13"; Value="0.09061"},
    @{Cell="E13"; Value="-0.74%"},
    @{Cell="D14"; Value="0.001539"},
    @{Cell="E14"; Value="0.66%"},
    @{Cell="D15"; Value="0.0006148"},
    @{Cell="E15"; Value="1.18%"},
    @{Cell="D16"; Value="0.006088"},
    @{Cell="E16"; Value="-0.09%"},
    @{Cell="D17"; Value="3.468"},
    @{Cell="E17"; Value="0.08%"},
    @{Cell="E18"; Value="-0.01%"},
    @{Cell="D19"; Value="2.178"},
    @{Cell="E19"; Value="0.55%"},
    @{Cell="E20"; Value="0.44%"},
    @{Cell="D21"; Value="0.1300"},
    @{Cell="E21"; Value="0.16%"},
    @{Cell="D22"; Value="4.090"},
    @{Cell="E22"; Value="-0.80%"},
    @{Cell="D23"; Value="0.04225"},
    @{Cell="E23"; Value="-0.32%"},
    @{Cell="D24"; Value="0.001213"},
    @{Cell="E24"; Value="-0.26%"},
    @{Cell="D25"; Value="0.004132"},
    @{Cell="E25"; Value="2.35%"},
    @{Cell="E26"; Value="0.08%"},
    @{Cell="E27"; Value="5.09%"},
    @{Cell="E40"; Value="1.33%"},
    @{Cell="D41"; Value="0.1112"},
    @{Cell="E41"; Value="-0.22%"},
    @{Cell="D42"; Value="0.004132"},
    @{Cell="E42"; Value="2.89%"},
    @{Cell="D43"; Value="0.002183"},
    @{Cell="E43"; Value="-0.73%"},
    @{Cell="D44"; Value="0.01351"},
    @{Cell="E44"; Value="-11.71%"},
    @{Cell="D45"; Value="0.00005175"},
    @{Cell="E45"; Value="1.44%"},
    @{Cell="E46"; Value="0.08%"},
    @{Cell="D47"; Value="0.03592"},
    @{Cell="E47"; Value="-34.14%"},
    @{Cell="E49"; Value="0.08%"},
    @{Cell="E50"; Value="0.08%"}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Formula = '="' + $u.Value + '"'
}

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.Copy()
    $c.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0

